# Apply the edits described by the diff:
# 1) Fix capitalization of the name in E2: "erit gridnev" -> "Erit Gridnev"
# 2) Remove the "OREN LAVI" time entry row (original row 4), shifting the
#    subsequent rows (5,6,7) up so they become rows 4,5,6. This also
#    shrinks the used range from A1:F7 to A1:F6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the name capitalization on row 2
$ws.Range("E2").Value = "Erit Gridnev"

# Delete the entire row 4 (OREN LAVI entry), shifting rows below it up
$ws.Rows.Item(4).Delete()
